# ForeignDD.xlsx update — add CREDIT.ACCT.NO and CHEQUE.NUMBER columns,
# switch the credit currency from EUR to USD, and populate the new fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B; everything that used to be in
# B:G (INSTRUMENT.TYPE:1 .. BK.TO.BK.OUT:1:1) shifts right to C:H.
[void]$ws.Columns.Item(2).Insert()

# Credit currency switches from EUR to USD.
$ws.Range("A2").Value = "USD"

# New column B: CREDIT.ACCT.NO header + the account number value.
$ws.Range("B1").Value = "CREDIT.ACCT.NO"
$ws.Range("B2").Value = 1000075724

# Column width for the new column B (closest the engine's column-width
# rounding allows to column A's 17.28515625 width).
$ws.Columns.Item(2).ColumnWidth = 16.45

# New trailing column I: CHEQUE.NUMBER header + the cheque number value.
$ws.Range("I1").Value = "CHEQUE.NUMBER"
$ws.Range("I2").Value = 123456

# Selection moves to the newly added account-number cell.
[void]$ws.Range("B2").Select()
